# Tuntikirjanpito.xlsx update
# - adds three new time-tracking entries (rows 120-122)
# - extends the "tunnit yht." SUM formula to cover the new rows
# - moves the sheet's scroll position / selection to the new bottom rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 120: reuse the date formatting (style) already used by A117 -------
$ws.Range("A117").Copy()
$ws.Range("A120").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A120").Value = 44596
$ws.Range("B120").Value = 2
$ws.Range("C120").Value = "NoteAddEditInput testit, NoteAddEditInputCntr luotu API/store logiikan erotteluun, komponentin mockauksen opettelua"
$ws.Range("D120").Value = "client"

# --- Row 121 -----------------------------------------------------------------
$ws.Range("B121").Value = 3
$ws.Range("C121").Value = "FavoritesList, FavoriteListItem testit"
$ws.Range("D121").Value = "client"

# --- Row 122 -----------------------------------------------------------------
$ws.Range("B122").Value = 1
$ws.Range("C122").Value = "Etusivun viimeisiä kommentteja ennen palautusta"
$ws.Range("D122").Value = "client"

# --- Update the totals formula to include the new rows -----------------------
$ws.Range("B123").Formula = "=SUM(B2:B122)"

# --- Move the view: scroll down a bit and select the last new row -----------
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 94
$win.ScrollColumn = 1
$ws.Range("C122").Select() | Out-Null
